$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update calibration parameter values (B and C columns) and E column sample counts
$ws.Range("B3").Value = 0.1
$ws.Range("C3").Value = -0.35
$ws.Range("E3").Value = 4

$ws.Range("B4").Value = 0.1
$ws.Range("C4").Value = -0.35
$ws.Range("E4").Value = 6

$ws.Range("B5").Value = 0.1
$ws.Range("C5").Value = -0.35

$ws.Range("B6").Value = 0.1
$ws.Range("C6").Value = -0.35
$ws.Range("E6").Value = 8

$ws.Range("B7").Value = 0.1
$ws.Range("C7").Value = -0.35
$ws.Range("E7").Value = 9

$ws.Range("B8").Value = 0.1
$ws.Range("C8").Value = -0.35
$ws.Range("E8").Value = 10

# Update the active cell selection on the sheet
$ws.Range("F14").Select()
